$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 "Save" - same style as other header cells (s="1")
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Data column H2:H6
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
